$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'24.793.91"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "'1.682.57"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.76%  "
$ws.Range("E4").Value = "  +0.63%  "
$ws.Range("D5").Value = "'313.78"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.39%  "
$ws.Range("D6").Value = "'1.004"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("D7").Value = "'0.3942"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("D8").Value = "'0.3970"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.46%  "
$ws.Range("D9").Value = "'1.007"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.80%  "
$ws.Range("D10").Value = "'1.421"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.37%  "
$ws.Range("D11").Value = "'51.69"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -3.48%  "
$ws.Range("D12").Value = "'0.08671"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.96%  "
$ws.Range("D13").Value = "'25.16"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.95%  "
$ws.Range("D14").Value = "'7.321"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.87%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.00001322"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.08%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'7.790"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -4.51%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "'94.15"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -3.30%  "
$ws.Range("B18").Value = "TRON"
$ws.Range("C18").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D18").Value = "'0.07107"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.55%  "
$ws.Range("B19").Value = "Avalanche"
$ws.Range("C19").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D19").Value = "'20.18"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.92%  "
$ws.Range("B20").Value = "WrappedEther"
$ws.Range("C20").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D20").Value = "'1.426.72"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -17.06%  "
$ws.Range("D21").Value = "'7.130"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.72%  "
$ws.Range("D22").Value = "'1.005"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.47%  "
$ws.Range("E23").Value = "  -2.46%  "
$ws.Range("D24").Value = "'24.784.45"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.19%  "
$ws.Range("D25").Value = "'2.360"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.98%  "
$ws.Range("D26").Value = "'23.81"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.44%  "
$ws.Range("D27").Value = "'2.771"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -8.33%  "
$ws.Range("D28").Value = "'162.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.45%  "
$ws.Range("D29").Value = "'150.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.06%  "
$ws.Range("D30").Value = "'5.735"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.61%  "
$ws.Range("D31").Value = "'2.585"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +16.23%  "
$ws.Range("D32").Value = "'7.834"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -8.64%  "
$ws.Range("B33").Value = "Hedera"
$ws.Range("C33").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D33").Value = "'0.08446"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.81%  "
$ws.Range("B34").Value = "VeChain"
$ws.Range("C34").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D34").Value = "'0.03074"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.96%  "
$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "'1.012"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -4.67%  "
$ws.Range("B36").Value = "InternetComputer(DFINITY)"
$ws.Range("C36").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D36").Value = "'6.943"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.50%  "
$ws.Range("B37").Value = "Algorand"
$ws.Range("C37").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D37").Value = "'0.2807"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.03%  "
$ws.Range("B38").Value = "Stellar"
$ws.Range("C38").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D38").Value = "'0.09576"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.40%  "
$ws.Range("B39").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C39").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D39").Value = "'1.580.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -16.96%  "
$ws.Range("D40").Value = "'10.45"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.19%  "
$ws.Range("D41").Value = "'0.7967"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -6.63%  "
$ws.Range("D42").Value = "'1.473"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("D43").Value = "'13.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.55%  "
$ws.Range("D44").Value = "'16.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.78%  "
$ws.Range("D45").Value = "'0.7168"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -4.54%  "
$ws.Range("D46").Value = "'2.582"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.68%  "
$ws.Range("D47").Value = "'4.184"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.32%  "
$ws.Range("D48").Value = "'0.08712"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.86%  "
$ws.Range("E49").Value = "  +0.40%  "
$ws.Range("D50").Value = "'1.343"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.89%  "
$ws.Range("D51").Value = "'138.06"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.96%  "
